$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.157908320426941
$ws.Range("B1").Value = 3.253323078155518
$ws.Range("C1").Value = 4.395672798156738
$ws.Range("D1").Value = 0.9752137660980225
$ws.Range("E1").Value = 1.190178394317627
